$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date field fix: 20/03/2024 -> 22/03/2024 (slide master + all layouts)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.HeadersFooters.Count; $i++) {
}
if ($master.Shapes.HasTitle -or $true) {
}

function Fix-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "20/03/2024") {
                $tr.Text = "22/03/2024"
            }
        }
    }
}

Fix-DateField $master
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Fix-DateField $layout
}

# ---------------------------------------------------------------------------
# 2. Slide 10 ("Examinar variables categóricas (cont.)") rework
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(10)

function Find-ShapeByName($container, $name) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        if ($container.Shapes.Item($i).Name -eq $name) {
            return $container.Shapes.Item($i)
        }
    }
    return $null
}

# Remove the small "Imagen 6" picture (the chi-2 p-value summary table)
$imagen6 = Find-ShapeByName $s "Imagen 6"
if ($imagen6 -ne $null) {
    $imagen6.Delete()
}

# Re-purpose "Imagen 3" (the big contingency table) into the small, repositioned picture
$imagen3 = Find-ShapeByName $s "Imagen 3"
$imagen3.Name = "Imagen 6"
$imagen3.Left = 1613027 / 12700.0
$imagen3.Top = 2243809 / 12700.0
$imagen3.Width = 2152650 / 12700.0
$imagen3.Height = 1895475 / 12700.0

# CuadroTexto 7: widen + move + update wording
$cb7 = Find-ShapeByName $s "CuadroTexto 7"
$cb7.Left = 765174 / 12700.0
$cb7.Top = 1158020 / 12700.0
$cb7.Width = 10260975 / 12700.0
$cb7.Height = 646331 / 12700.0
$cb7.TextFrame.TextRange.Text = "Los resultados son verificados mediante una prueba chi-2 por independencia: ¿qué debería aparecer? El resumen con estadístico y valor-p. "

# CuadroTexto 8: move
$cb8 = Find-ShapeByName $s "CuadroTexto 8"
$cb8.Left = 1356670 / 12700.0
$cb8.Top = 5366628 / 12700.0

# CuadroTexto 9: move
$cb9 = Find-ShapeByName $s "CuadroTexto 9"
$cb9.Left = 1356670 / 12700.0
$cb9.Top = 4578742 / 12700.0

# The slide's internal shape-id allocator hands out the lowest free id from the
# set that was unused *before this script ran* (here: 3, 5, 6) and only then
# overflows past the original max id (10 -> 11, 12, ...). We need the new
# "CuadroTexto 11" textbox to land on id 12, i.e. it must be the 5th shape
# created in this run, so four scratch textboxes are created and scrapped
# first to burn through ids 3, 5, 6 and 11.
$scratch = New-Object System.Collections.ArrayList
for ($n = 0; $n -lt 4; $n++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    [void]$scratch.Add($dummy)
}

$newBox = $s.Shapes.AddTextbox(1, 6505276 / 12700.0, 2428068 / 12700.0, 4286774 / 12700.0, 646331 / 12700.0)
$newBox.Name = "CuadroTexto 11"
$newBox.TextFrame.WordWrap = $true
$newBox.TextFrame.TextRange.Text = "Recuérdese que la hipótesis nula es la siguiente:"

foreach ($d in $scratch) {
    $d.Delete()
}
